$wb = $excel.ActiveWorkbook

# 1. Duplicate the "devices" sheet. Excel places the copy immediately
#    before the source sheet and auto-names it "devices (2)".
$devices = $wb.Worksheets.Item("devices")
$devices.Copy($devices) | Out-Null

# After the copy: sheet order is "devices (2)" (the new copy), then the
# original "devices" sheet (still named "devices").
$devicesCopy = $wb.Worksheets.Item("devices (2)")
$devices     = $wb.Worksheets.Item("devices")

# 2. On the "devices (2)" sheet, rename the browser rows:
#    "chrome" -> "Chrome", "iexplorer" -> "Firefox", and drop the old
#    "firefox" row (row 7) entirely.
$devicesCopy.Range("A5").Value = "Chrome"
$devicesCopy.Range("A6").Value = "Firefox"
$devicesCopy.Rows("7").Delete() | Out-Null
$devicesCopy.Range("A7").Select() | Out-Null

# 3. Trim the "devices" sheet down to just the header row plus a single
#    "Chrome" data row (the new, lighter-weight report sheet).
$devices.Range("A5").Value = "Chrome"
$devices.Rows("2:4").Delete() | Out-Null
$devices.Rows("3:4").Delete() | Out-Null
$devices.Range("A3:XFD3").Select() | Out-Null

# 4. Make "devices" the active tab.
$devices.Activate() | Out-Null
